$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.744.69'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -3.59%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.487.25'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -6.19%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '556.33'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -4.44%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.46'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -5.36%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.601'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -3.46%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.486.67'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -6.08%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -8.17%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.51'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -5.36%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.359'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -6.41%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.57'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -7.43%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -6.05%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -8.71%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.635.77'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.64%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.485.37'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -5.93%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -7.75%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.22'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -7.20%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.25'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -6.70%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '323.53'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -6.75%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +3.17%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '64.43'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -5.36%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -9.40%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'Bittensor'
$ws.Range('B27').Style = "Normal"
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C27').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '559.30'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -5.71%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('B28').Style = "Normal"
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('C28').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.603.11'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -6.18%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.52'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -6.04%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.96'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -3.41%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.43'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -9.78%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.46%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -7.18%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.61'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -7.88%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.99'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -10.34%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.96'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -10.29%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.385'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -5.61%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '147.11'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -7.70%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '40.61'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -3.16%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -5.01%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '149.06'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -9.08%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('B47').Style = "Normal"
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C47').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.68'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -6.19%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('B48').Style = "Normal"
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C48').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.20'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -8.83%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0546'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -7.79%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -5.62%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -5.34%  '
$ws.Range('E51').Style = "Normal"
